$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Drop the leftover "_GoBack" bookmark in the first (empty) paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Merge runs that were only split because of tracked-change history.
#    Re-finding the full (already contiguous) sentence and "replacing" it
#    with itself collapses the surrounding runs into a single run, which is
#    exactly what the target markup looks like.
# ---------------------------------------------------------------------------
$mergeTexts = @(
    "This course is intended to teach the necessary knowledge and skills to install, configure, and use the IBM Aspera High-Speed Transfer Server.",
    "Classroom or self-paced virtual classroom (SPVC)",
    "After completing this course, you should be able to:",
    "Perform file transfers using the Aspera GUI and from the command line",
    "Fundamental knowledge of using Windows and Linux operating systems",
    "This unit describes the operation of the FASP protocol and how it compares with traditional file transfer protocols.",
    "Clarify the value of using Vlinks",
    "This unit a brief overview of IBM Aspera software and how they may be integrated as an Aspera environment",
    "Locate and install the appropriate IBM Aspera High-Speed Transfer Server software (Windows and Linux)",
    "Verify installation success by transferring files to and from the Aspera Demo Server",
    "This exercise demonstrates the installation of IBM Aspera High-Speed Transfer (HST) Server software.",
    "Unit 4. Configuring IBM Aspera High-Speed Transfer Server",
    "Manage file permissions for inbound/outbound transfers",
    "Create predefined connections and share them with other users",
    "This unit addresses the basic configuration for adding and managing Aspera transfer users and groups",
    "Verify user account’s ability to perform FASP-based transfers",
    "This exercise uses the Aspera GUI to configure new users and groups, to define configuration parameters that manage transfers performed by users and groups, and to implement Vlinks.",
    "Unit 6. Using command-line operations",
    "Transfer files and directories between Aspera servers using the ascp command",
    "This exercise demonstrates the use of IBM Aspera command-line tools to configure the aspera.conf file with the asconfigurator utility and to initiate and manage file transfers using the ascp command.",
    "Unit 7. Configuring advanced features",
    "This unit addresses several features that are not required for basic configuration of the IBM Aspera Transfer Server, but are commonly implemented on production systems",
    "Configure IBM Aspera High-Speed Transfer Server to use custom SSL certificates and token authorization",
    "Explain the procedure for implementing hot folders on Windows platforms",
    "Unit 8. Routine maintenance tasks",
    "This unit identifies common performance bottlenecks, presents common maintenance tasks, and introduces how to interpret some of the Aspera log file entries."
)

foreach ($t in $mergeTexts) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2)
    if (-not $ok) {
        Write-Output "WARNING: merge text not found: $t"
    }
}

# ---------------------------------------------------------------------------
# 3. Remove the "stay informed" social-media block from the "For more
#    information" section, leaving a single empty (indented) paragraph in
#    its place.
# ---------------------------------------------------------------------------
$rStart = $d.Content
$rStart.Find.Execute("To stay informed about IBM training") | Out-Null
$pFirst = $rStart.Paragraphs(1)

$rEnd = $d.Content
$rEnd.Find.Execute("twitter.com/websphere_edu") | Out-Null
$pLast = $rEnd.Paragraphs(1)

# Delete the four leading paragraphs entirely (their text and paragraph
# marks), stopping right at the start of the final ("Twitter: ...") one.
$dropRange = $d.Range($pFirst.Range.Start, $pLast.Range.Start)
$dropRange.Delete()

# Clear the remaining paragraph's text, keeping its paragraph mark/pPr
# (pStyle Abstractbodytext, keepNext, keepLines, ind left=720) intact.
$rLast = $d.Content
$rLast.Find.Execute("twitter.com/websphere_edu") | Out-Null
$pLast2 = $rLast.Paragraphs(1)
$textRange = $d.Range($pLast2.Range.Start, $pLast2.Range.End - 1)
$textRange.Delete()

Write-Output "done"
